# Refresh the cryptocurrency price / 1h-volume snapshot (GitHub Actions data pull).
# Row 29/30 (FirstDigitalUSD / NEARProtocol) also swap rank position in this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.791.12'
$ws.Range('E2').Value = '  -1.41%  '
$ws.Range('D3').Value = '3.031.44'
$ws.Range('E3').Value = '  -1.82%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.63'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.71'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -4.18%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('E8').Value = '  -3.09%  '
$ws.Range('D9').Value = '3.031.75'
$ws.Range('E9').Value = '  -1.77%  '
$ws.Range('E10').Value = '  -3.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.68'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.82%  '
$ws.Range('E12').Value = '  -2.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000232'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -3.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.47'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -5.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.120'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.23%  '
$ws.Range('D16').Value = '3.533.94'
$ws.Range('E16').Value = '  -1.89%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.08'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.29%  '
$ws.Range('D18').Value = '62.750.47'
$ws.Range('E18').Value = '  -1.49%  '
$ws.Range('D19').Value = '3.031.43'
$ws.Range('E19').Value = '  -1.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '468.59'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.08'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.693'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.61%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.42'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.40'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.14'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.44'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.36%  '
$ws.Range('E27').Value = '  +2.30%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('B29').Value = 'NEARProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.27'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.45%  '
$ws.Range('B30').Value = 'FirstDigitalUSD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('E31').Value = '  -1.89%  '
$ws.Range('E32').Value = '  -1.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.50'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.73%  '
$ws.Range('E34').Value = '  -4.55%  '
$ws.Range('E35').Value = '  -0.95%  '
$ws.Range('D36').Value = '0.0₃0800'
$ws.Range('E36').Value = '  -5.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.79'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -4.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.16'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '50.31'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.98%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.99'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -13.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.03'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '423.73'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -5.15%  '
$ws.Range('E43').Value = '  -1.70%  '
$ws.Range('E44').Value = '  +0.85%  '
$ws.Range('D45').Value = '2.801.42'
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0357'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '37.89'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -10.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.93'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.13%  '
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.57'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.59%  '
$ws.Range('E51').Value = '  -1.52%  '
